# Auto-generated Excel COM-interop script
# Updates the "想去人数" (F column) values on every worksheet to match
# the commit "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 211
$ws.Range("F3").Value = 694
$ws.Range("F4").Value = 540
$ws.Range("F5").Value = 2217
$ws.Range("F6").Value = 1316
$ws.Range("F7").Value = 787
$ws.Range("F8").Value = 83
$ws.Range("F10").Value = 2836
$ws.Range("F11").Value = 21
$ws.Range("F13").Value = 1070
$ws.Range("F14").Value = 571
$ws.Range("F15").Value = 520
$ws.Range("F16").Value = 899
$ws.Range("F17").Value = 89
$ws.Range("F18").Value = 91
$ws.Range("F19").Value = 532
$ws.Range("F20").Value = 112
$ws.Range("F21").Value = 603
$ws.Range("F23").Value = 268
$ws.Range("F25").Value = 960
$ws.Range("F26").Value = 4872
$ws.Range("F27").Value = 376
$ws.Range("F28").Value = 138
$ws.Range("F29").Value = 64

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 14
$ws.Range("F3").Value = 72
$ws.Range("F6").Value = 353
$ws.Range("F7").Value = 384
$ws.Range("F14").Value = 281
$ws.Range("F21").Value = 26
$ws.Range("F23").Value = 34
$ws.Range("F24").Value = 34
$ws.Range("F25").Value = 345
$ws.Range("F27").Value = 547
$ws.Range("F37").Value = 701

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value = 373
$ws.Range("F7").Value = 343

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 211
$ws.Range("F6").Value = 373
$ws.Range("F7").Value = 14
$ws.Range("F8").Value = 695
$ws.Range("F9").Value = 353
$ws.Range("F10").Value = 384
$ws.Range("F11").Value = 540
$ws.Range("F12").Value = 2217
$ws.Range("F13").Value = 1316
$ws.Range("F14").Value = 787
$ws.Range("F15").Value = 83
$ws.Range("F19").Value = 2837
$ws.Range("F20").Value = 21
$ws.Range("F21").Value = 281
$ws.Range("F23").Value = 1070
$ws.Range("F24").Value = 571
$ws.Range("F26").Value = 343
$ws.Range("F27").Value = 520
$ws.Range("F28").Value = 899
$ws.Range("F29").Value = 899
$ws.Range("F30").Value = 89
$ws.Range("F31").Value = 26
$ws.Range("F33").Value = 91
$ws.Range("F34").Value = 112
$ws.Range("F35").Value = 34
$ws.Range("F36").Value = 34
$ws.Range("F37").Value = 603
$ws.Range("F39").Value = 345
$ws.Range("F40").Value = 547
$ws.Range("F41").Value = 268
$ws.Range("F44").Value = 960
$ws.Range("F45").Value = 4872
$ws.Range("F47").Value = 376
$ws.Range("F48").Value = 138
$ws.Range("F49").Value = 701
